# Commit: "a minor data table structure change"
# Move CoopNodeID and PickNodeID columns from LineItem table to Skus Table
#
# In this workbook each worksheet models one Route as a single column of
# node names. The "move columns" restructuring shows up here as each
# sheet's route column being extended with the additional
# CoopNode/PickNode/Input@Depot/EndNode steps that used to live in a
# separate table. We simply append the extra rows under the existing
# "Route" column on every sheet.

$wb = $excel.ActiveWorkbook

# --- Picker1 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Picker1")
$ws.Range("A10").Value = "MyPickNode11"
$ws.Range("A11").Value = "MyCoopNode6"
$ws.Range("A12").Value = "MyPickNode11"
$ws.Range("A13").Value = "MyCoopNode6"
$ws.Range("A14").Value = "MyPickNode12"
$ws.Range("A15").Value = "MyCoopNode6"
$ws.Range("A16").Value = "MyPickNode3"
$ws.Range("A17").Value = "MyCoopNode2"
$ws.Range("A18").Value = "MyPickNode9"
$ws.Range("A19").Value = "MyCoopNode5"
$ws.Range("A20").Value = "MyPickNode13"
$ws.Range("A21").Value = "MyCoopNode7"
$ws.Range("A22").Value = "MyPickNode19"
$ws.Range("A23").Value = "MyCoopNode10"
$ws.Range("A24").Value = "EndNode"
$ws.Range("A14").Select()

# --- Picker2 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Picker2")
$ws.Range("A8").Value = "MyPickNode23"
$ws.Range("A9").Value = "MyCoopNode12"
$ws.Range("A10").Value = "MyPickNode22"
$ws.Range("A11").Value = "MyCoopNode11"
$ws.Range("A12").Value = "MyPickNode21"
$ws.Range("A13").Value = "MyCoopNode11"
$ws.Range("A14").Value = "MyPickNode21"
$ws.Range("A15").Value = "MyCoopNode11"
$ws.Range("A16").Value = "MyPickNode14"
$ws.Range("A17").Value = "MyCoopNode7"
$ws.Range("A18").Value = "MyPickNode14"
$ws.Range("A19").Value = "MyCoopNode7"
$ws.Range("A20").Value = "EndNode"
$ws.Rows.Item(7).RowHeight = 12.75
$ws.Rows("8:9").Select()

# --- Transporter1 ------------------------------------------------------
$ws = $wb.Worksheets.Item("Transporter1")
$ws.Range("A4").Value = "MyCoopNode6"
$ws.Range("A5").Value = "Input@Depot"
$ws.Range("A6").Value = "MyCoopNode7"
$ws.Range("A7").Value = "Input@Depot"
$ws.Range("A8").Value = "EndNode"
$ws.Range("K18").Select()

# --- Transporter2 ------------------------------------------------------
$ws = $wb.Worksheets.Item("Transporter2")
$ws.Range("A5").Value = "MyCoopNode2"
$ws.Range("A6").Value = "MyCoopNode5"
$ws.Range("A7").Value = "MyCoopNode7"
$ws.Range("A8").Value = "Input@Depot"
$ws.Range("A9").Value = "EndNode"
$ws.Range("A9").Select()

# --- Transporter3 (ends as the active / tab-selected sheet) ------------
$ws = $wb.Worksheets.Item("Transporter3")
$ws.Range("A5").Value = "MyCoopNode12"
$ws.Range("A6").Value = "MyCoopNode11"
$ws.Range("A7").Value = "Input@Depot"
$ws.Range("A8").Value = "MyCoopNode10"
$ws.Range("A9").Value = "Input@Depot"
$ws.Range("A10").Value = "EndNode"
$ws.Activate()
$ws.Range("G17").Select()
